$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 11:32"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 5415977
$ws.Range("C4").Value = 311
$ws.Range("D4").Value = 2843642
$ws.Range("E4").Value = 2401916
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 170419

# Row 26: Indonesia -> Indonesia
$ws.Range("B26").Value = 135123
$ws.Range("C26").Value = 2307
$ws.Range("D26").Value = 89618
$ws.Range("E26").Value = 39484
$ws.Range("G26").Value = 53
$ws.Range("H26").Value = 6021

# Row 33: Israel -> Israel
$ws.Range("B33").Value = 90472
$ws.Range("C33").Value = 650
$ws.Range("D33").Value = 66151
$ws.Range("E33").Value = 23664
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 657

# Row 48: Polonia -> Polonia
$ws.Range("B48").Value = 55312
$ws.Range("C48").Value = 825
$ws.Range("D48").Value = 38362
$ws.Range("E48").Value = 15092
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 1858

# Row 58: Afganistan -> Afganistan
$ws.Range("B58").Value = 37431
$ws.Range("C58").Value = 7
$ws.Range("E58").Value = 9354

# Row 65: Serbia -> Serbia
$ws.Range("D65").Value = 26117
$ws.Range("E65").Value = 2220

# Row 71: Australia -> Austria
$ws.Range("A71").Value = "Austria"
$ws.Range("B71").Value = 22876
$ws.Range("C71").Value = 282
$ws.Range("D71").Value = 20499
$ws.Range("E71").Value = 1652
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 725

# Row 72: Austria -> Australia
$ws.Range("A72").Value = "Australia"
$ws.Range("B72").Value = 22743
$ws.Range("C72").Value = 385
$ws.Range("D72").Value = 13350
$ws.Range("E72").Value = 9018
$ws.Range("G72").Value = 14
$ws.Range("H72").Value = 375

# Row 73: El Salvador -> El Salvador
$ws.Range("B73").Value = 22314
$ws.Range("C73").Value = 321
$ws.Range("D73").Value = 10455
$ws.Range("E73").Value = 11264
$ws.Range("G73").Value = 11
$ws.Range("H73").Value = 595

# Row 87: Consejo Danes para los Refugiados -> Consejo Danes para los Refugiados
$ws.Range("B87").Value = 9605
$ws.Range("C87").Value = 16
$ws.Range("D87").Value = 8512
$ws.Range("E87").Value = 855
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 238

# Row 88: Malasia -> Malasia
$ws.Range("B88").Value = 9149
$ws.Range("C88").Value = 20
$ws.Range("D88").Value = 8828
$ws.Range("E88").Value = 196

# Row 97: Finlandia -> Finlandia
$ws.Range("B97").Value = 7700
$ws.Range("C97").Value = 17
$ws.Range("E97").Value = 317

# Row 111: Hong Kong -> Hong Kong
$ws.Range("B111").Value = 4361
$ws.Range("C111").Value = 48
$ws.Range("D111").Value = 3392
$ws.Range("E111").Value = 903

# Row 128: Eslovenia -> Eslovenia
$ws.Range("B128").Value = 2369
$ws.Range("C128").Value = 37
$ws.Range("E128").Value = 280

# Row 129: Lituania -> Lituania
$ws.Range("B129").Value = 2352
$ws.Range("C129").Value = 22
$ws.Range("D129").Value = 1691
$ws.Range("E129").Value = 580

# Row 131: Estonia -> Estonia
$ws.Range("B131").Value = 2177
$ws.Range("C131").Value = 3
$ws.Range("D131").Value = 1976
$ws.Range("E131").Value = 138

# Row 213: Montserrat -> Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214: Islas Malvinas -> Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

